$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '56.212.95'
$ws.Range("E2").Value = '  +9.21%  '
$ws.Range("D3").Value = '3.231.66'
$ws.Range("E3").Value = '  +4.10%  '
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '398.87'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.96%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '110.80'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +7.53%  '
$ws.Range("E8").Value = '  +0.00%  '
$ws.Range("E9").Value = '  +7.57%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '39.63'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +6.59%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0896'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +5.30%  '
$ws.Range("E12").Value = '  +2.05%  '
$ws.Range("D13").Value = '3.750.43'
$ws.Range("E13").Value = '  +4.35%  '
$ws.Range("E14").Value = '  +2.87%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '8.09'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +3.62%  '
$ws.Range("B16").Value = 'WrappedEther'
$ws.Range("C16").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D16").Value = '3.241.49'
$ws.Range("E16").Value = '  +4.48%  '
$ws.Range("B17").Value = 'Polygon'
$ws.Range("C17").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.07'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +7.47%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '10.55'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -6.52%  '
$ws.Range("D19").Value = '56.101.37'
$ws.Range("E20").Value = '  +2.35%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.12'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +6.66%  '
$ws.Range("E22").Value = '  +4.89%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '289.42'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +8.93%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '74.13'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +6.00%  '
$ws.Range("E25").Value = '  +4.64%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.19'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.57%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '28.31'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +4.92%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.54'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +4.03%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.172'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +3.22%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.999'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.05%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.111'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +4.34%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '11.28'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +9.38%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0498'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +5.87%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '36.99'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +4.40%  '
$ws.Range("E35").Value = '  +1.76%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '51.06'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.49%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.58'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +6.88%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.00'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.07%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.08'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +21.86%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '139.39'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +8.10%  '
$ws.Range("E41").Value = '  +2.43%  '
$ws.Range("E42").Value = '  +10.11%  '
$ws.Range("E43").Value = '  -3.55%  '
$ws.Range("E44").Value = '  +2.15%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '16.78'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.66%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '22.66'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.65%  '
$ws.Range("B47").Value = 'ThetaToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.11'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +42.37%  '
$ws.Range("B48").Value = 'ApeXProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.47'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.36%  '
$ws.Range("B49").Value = 'WEMIXToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.09'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.11%  '
$ws.Range("D50").Value = '2.133.04'
$ws.Range("E50").Value = '  +3.65%  '
$ws.Range("B51").Value = 'BEAM'
$ws.Range("C51").Value = 'https://coinranking.com/coin/cYYMfXF4u+beam-beam'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0343'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +5.67%  '
